$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# New K-Prim example block, rows 68-76 (mirrors the existing SC example at
# rows 39-51). Values first, then formats copied from the matching template
# cells so the new rows reuse the same style entries as the rest of the sheet.
# ---------------------------------------------------------------------------

# Row 68 - section header "Typ" / "KPRIM" / "Fragetpy: K-Prim"
$ws.Range("A68").Value = "Typ"
$ws.Range("B68").Value = "KPRIM"
$ws.Range("D68").Value = "Fragetpy: K-Prim"

# Row 69 - Title
$ws.Range("A69").Value = "Title"
$ws.Range("B69").Value = "Fussball: Weltmeister"
$ws.Range("D69").Value = "Wird in Navigation angezeigt. Entspricht dem Thema der Frage. "

# Row 70 - Description
$ws.Range("A70").Value = "Description"
$ws.Range("B70").Value = "Prüfen Sie die Weltmeiser kennen"
$ws.Range("D70").Value = "Optionale Beschreibung / Einführungstext der zwischen dem Titel und der eigentlichen Frage erscheint. Worum geht es in dieser Frage?"

# Row 71 - Question
$ws.Range("A71").Value = "Question"
$ws.Range("B71").Value = "Die folgenden Länder haben die Fussball Weltmeistertitel bereits mehr als einmal gewonnen."
$ws.Range("D71").Value = "Die Frage"

# Row 72 - Points
$ws.Range("A72").Value = "Points"
$ws.Range("B72").Value = 1
$ws.Range("D72").Value = "Die maximal zu erreichende Punktzahl"

# Row 73 - true answer "Deutschland"
$ws.Range("A73").Value = "+"
$ws.Range("B73").Value = "Deutschland"
$ws.Range("D73").Value = "Eine wahre Antwort"

# Row 74 - false answer "Frankreich"
$ws.Range("A74").Value = "-"
$ws.Range("B74").Value = "Frankreich"
$ws.Range("D74").Value = "Eine falsche Antwort"

# Row 75 - false answer "Spanien"
$ws.Range("A75").Value = "-"
$ws.Range("B75").Value = "Spanien"
$ws.Range("D75").Value = "Eine falsche Antwort"

# Row 76 - true answer "Uruguay"
$ws.Range("A76").Value = "+"
$ws.Range("B76").Value = "Uruguay"
$ws.Range("D76").Value = "Eine wahre Antwort"

# ---------------------------------------------------------------------------
# Copy cell formatting from the analogous rows of the existing SC example
# (rows 39-51) so the new block's cellXfs line up with the rest of the sheet.
# ---------------------------------------------------------------------------
function Copy-RowFormat($srcRow, $dstRow) {
    foreach ($col in @("A", "B", "D")) {
        $ws.Range("$col$srcRow").Copy() | Out-Null
        $ws.Range("$col$dstRow").PasteSpecial(-4122) | Out-Null
    }
}

Copy-RowFormat 39 68
Copy-RowFormat 40 69
Copy-RowFormat 41 70
Copy-RowFormat 42 71
Copy-RowFormat 43 72
Copy-RowFormat 44 73
Copy-RowFormat 47 74
Copy-RowFormat 47 75
Copy-RowFormat 44 76

$excel.CutCopyMode = 0

# Row heights for the wrapped description rows (matches rows 41/42 in the
# template block).
$ws.Rows.Item(70).RowHeight = 30
$ws.Rows.Item(71).RowHeight = 45

# ---------------------------------------------------------------------------
# View/selection bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("B81").Select()
$excel.ActiveWindow.ScrollRow = 23
